# Weekly data refresh: insert one new daily price record for Berenjena
# (Vega Central Mapocho de Santiago) at row 205, pushing the existing
# records at rows 205-255 down to 206-256.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 205; Excel shifts rows 205:255 -> 206:256
# and grows the sheet's used range (and <dimension>) to A1:R256 automatically.
$ws.Rows.Item(205).Insert()

# Populate the newly inserted row 205 with the new record's data.
$ws.Cells.Item(205, 1).Value  = 9
$ws.Cells.Item(205, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(205, 3).Value  = "Metropolitana"
$ws.Cells.Item(205, 4).Value  = 44641
$ws.Cells.Item(205, 5).Value  = 13
$ws.Cells.Item(205, 6).Value  = 100112001
$ws.Cells.Item(205, 7).Value  = "Berenjena"
$ws.Cells.Item(205, 8).Value  = "Sin especificar"
$ws.Cells.Item(205, 9).Value  = "Primera"
$ws.Cells.Item(205, 10).Value = 97
$ws.Cells.Item(205, 11).Value = 8000
$ws.Cells.Item(205, 12).Value = 9000
$ws.Cells.Item(205, 13).Value = 8495
$ws.Cells.Item(205, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(205, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(205, 16).Value = 170
$ws.Cells.Item(205, 17).Value = 50
$ws.Cells.Item(205, 18).Value = "Hortaliza"
